$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-15 and 18-19 already carry the "bordered + wrap/centered" style (style index 2)
# for columns B,E,F and just need values + a date format added for C,D.
$plainRows = @(9, 10, 11, 12, 13, 14, 15, 18, 19)

foreach ($r in $plainRows) {
    $ws.Range("B$r").Value = "Chentao Jin"

    $ws.Range("C$r").Value = 45386
    $ws.Range("C$r").NumberFormat = "d-mmm"

    $ws.Range("D$r").Value = 45387
    $ws.Range("D$r").NumberFormat = "d-mmm"

    $ws.Range("E$r").Value = 0.3
    $ws.Range("F$r").Value = 0.3
}

# Rows 20-21 currently use the "unbordered / no-alignment" style (style index 3) for B:F,
# so first bring their formatting in line with the rows above (copy from row 19, which
# already has the right border + wrap/centered + date formatting), then fill in the values.
$targetRows = @(20, 21)

foreach ($r in $targetRows) {
    $ws.Range("B19:F19").Copy()
    $ws.Range("B$r`:F$r").PasteSpecial(-4122)

    $ws.Range("B$r").Value = "Chentao Jin"

    $ws.Range("C$r").Value = 45386
    $ws.Range("D$r").Value = 45387

    $ws.Range("E$r").Value = 0.3
    $ws.Range("F$r").Value = 0.3
}

$excel.CutCopyMode = $false

[void]$ws.Range("G26").Select()

Write-Host "team plan sheet filled"
